$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1956
$ws.Range("J3").Value = 2039
$ws.Range("I4").Value = 1754
$ws.Range("J4").Value = 455
$ws.Range("J5").Value = 151
$ws.Range("J6").Value = 2554
$ws.Range("I7").Value = 26202
$ws.Range("J7").Value = 7155

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 52
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 18
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 185

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 204
$ws.Range("J8").Value = 444
$ws.Range("J9").Value = 47
$ws.Range("J10").Value = 42
$ws.Range("J11").Value = 93
$ws.Range("J18").Value = 84
$ws.Range("J19").Value = 243
$ws.Range("J20").Value = 150
$ws.Range("I22").Value = 68
$ws.Range("J27").Value = 42
$ws.Range("J29").Value = 398
$ws.Range("J31").Value = 48
$ws.Range("J32").Value = 13
$ws.Range("J33").Value = 300
$ws.Range("J34").Value = 42
$ws.Range("J36").Value = 107
$ws.Range("J37").Value = 239
$ws.Range("J42").Value = 279
$ws.Range("J47").Value = 65
$ws.Range("J48").Value = 65
$ws.Range("J53").Value = 66
$ws.Range("J55").Value = 83
$ws.Range("J63").Value = 27
$ws.Range("J65").Value = 185
$ws.Range("J67").Value = 260
$ws.Range("J68").Value = 15
$ws.Range("J71").Value = 31
$ws.Range("J76").Value = 107
$ws.Range("J78").Value = 92
$ws.Range("J79").Value = 218
$ws.Range("J84").Value = 74
$ws.Range("J85").Value = 333
$ws.Range("J86").Value = 40
$ws.Range("J89").Value = 77
$ws.Range("I91").Value = 278
$ws.Range("J95").Value = 106
$ws.Range("J96").Value = 80
$ws.Range("J99").Value = 93
$ws.Range("I101").Value = 26202
$ws.Range("J101").Value = 7155

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 77
$ws.Range("J3").Value = 87
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 117
$ws.Range("J3").Value = 141
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 398

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 61
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 243

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 135
$ws.Range("J7").Value = 333

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J5").Value = 6
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 61
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 22
$ws.Range("J3").Value = 11
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 34
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J5").Value = 12
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 444

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 204
